$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-CellXml($cell, [string]$innerParagraphXml) {
    $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $innerParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $cell.Range.InsertXML($xml)
}

# IDS / No row -> "# Paths Popped from Queue" cell: 3479 -> 4888
$cell = $tbl.Cell(6, 4)
$p = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t xml:space="preserve">           </w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>4888</w:t></w:r></w:p>'
Set-CellXml $cell $p

# IDS / Yes row -> "Time (s)" cell: 1.52 x 10^-4 -> 1.12 x 10^-3 (first number split into 3 runs)
$cell = $tbl.Cell(7, 3)
$p = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>1.</w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t xml:space="preserve"> x 10</w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>3</w:t></w:r></w:p>'
Set-CellXml $cell $p

# IDS / Yes row -> "# Paths Popped from Queue" cell: 57 -> 749 (leading spaces 13 -> 12)
$cell = $tbl.Cell(7, 4)
$p = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">            </w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>749</w:t></w:r></w:p>'
Set-CellXml $cell $p

# IDS / Yes row -> "Max Queue Size" cell: 30 -> 53 (paragraph rPr simplified)
$cell = $tbl.Cell(7, 5)
$p = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">        </w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>53</w:t></w:r></w:p>'
Set-CellXml $cell $p

# IDS / Yes row -> "Returned Path's Length/Cost" cell: 9 -> 11 (leading spaces 11 -> 10, paragraph rPr simplified)
$cell = $tbl.Cell(7, 6)
$p = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">          </w:t></w:r><w:r><w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t>11</w:t></w:r></w:p>'
Set-CellXml $cell $p

Write-Host "Done"
